$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45889.01041666666
$ws.Cells.Item(3, 1).Value = 45889.02083333334
$ws.Cells.Item(4, 1).Value = 45889.03125
$ws.Cells.Item(5, 1).Value = 45889.04166666666
$ws.Cells.Item(6, 1).Value = 45889.05208333334
$ws.Cells.Item(7, 1).Value = 45889.0625
$ws.Cells.Item(8, 1).Value = 45889.07291666666
$ws.Cells.Item(9, 1).Value = 45889.08333333334
$ws.Cells.Item(10, 1).Value = 45889.09375
$ws.Cells.Item(11, 1).Value = 45889.10416666666
$ws.Cells.Item(12, 1).Value = 45889.11458333334
$ws.Cells.Item(13, 1).Value = 45889.125
$ws.Cells.Item(14, 1).Value = 45889.13541666666
$ws.Cells.Item(15, 1).Value = 45889.14583333334
$ws.Cells.Item(16, 1).Value = 45889.15625
$ws.Cells.Item(17, 1).Value = 45889.16666666666
$ws.Cells.Item(18, 1).Value = 45889.17708333334
$ws.Cells.Item(18, 2).Value = 6
$ws.Cells.Item(19, 1).Value = 45889.1875
$ws.Cells.Item(19, 2).Value = 6
$ws.Cells.Item(20, 1).Value = 45889.19791666666
$ws.Cells.Item(20, 2).Value = 6
$ws.Cells.Item(21, 1).Value = 45889.20833333334
$ws.Cells.Item(21, 2).Value = 6
$ws.Cells.Item(22, 1).Value = 45889.21875
$ws.Cells.Item(22, 2).Value = 28
$ws.Cells.Item(23, 1).Value = 45889.22916666666
$ws.Cells.Item(23, 2).Value = 29
$ws.Cells.Item(24, 1).Value = 45889.23958333334
$ws.Cells.Item(24, 2).Value = 33
$ws.Cells.Item(25, 1).Value = 45889.25
$ws.Cells.Item(25, 2).Value = 40
$ws.Cells.Item(26, 1).Value = 45889.26041666666
$ws.Cells.Item(26, 2).Value = 202
$ws.Cells.Item(27, 1).Value = 45889.27083333334
$ws.Cells.Item(27, 2).Value = 224
$ws.Cells.Item(28, 1).Value = 45889.28125
$ws.Cells.Item(28, 2).Value = 253
$ws.Cells.Item(29, 1).Value = 45889.29166666666
$ws.Cells.Item(29, 2).Value = 289
$ws.Cells.Item(30, 1).Value = 45889.30208333334
$ws.Cells.Item(30, 2).Value = 696
$ws.Cells.Item(31, 1).Value = 45889.3125
$ws.Cells.Item(31, 2).Value = 732
$ws.Cells.Item(32, 1).Value = 45889.32291666666
$ws.Cells.Item(32, 2).Value = 812
$ws.Cells.Item(33, 1).Value = 45889.33333333334
$ws.Cells.Item(33, 2).Value = 868
$ws.Cells.Item(34, 1).Value = 45889.34375
$ws.Cells.Item(34, 2).Value = 1353
$ws.Cells.Item(35, 1).Value = 45889.35416666666
$ws.Cells.Item(35, 2).Value = 1439
$ws.Cells.Item(36, 1).Value = 45889.36458333334
$ws.Cells.Item(36, 2).Value = 1499
$ws.Cells.Item(37, 1).Value = 45889.375
$ws.Cells.Item(37, 2).Value = 1542
$ws.Cells.Item(38, 1).Value = 45889.38541666666
$ws.Cells.Item(38, 2).Value = 1851
$ws.Cells.Item(39, 1).Value = 45889.39583333334
$ws.Cells.Item(39, 2).Value = 1907
$ws.Cells.Item(40, 1).Value = 45889.40625
$ws.Cells.Item(40, 2).Value = 1946
$ws.Cells.Item(41, 1).Value = 45889.41666666666
$ws.Cells.Item(41, 2).Value = 1979
$ws.Cells.Item(42, 1).Value = 45889.42708333334
$ws.Cells.Item(42, 2).Value = 2166
$ws.Cells.Item(43, 1).Value = 45889.4375
$ws.Cells.Item(43, 2).Value = 2190
$ws.Cells.Item(44, 1).Value = 45889.44791666666
$ws.Cells.Item(44, 2).Value = 2207
$ws.Cells.Item(45, 1).Value = 45889.45833333334
$ws.Cells.Item(45, 2).Value = 2221
$ws.Cells.Item(46, 1).Value = 45889.46875
$ws.Cells.Item(46, 2).Value = 2292
$ws.Cells.Item(47, 1).Value = 45889.47916666666
$ws.Cells.Item(47, 2).Value = 2299
$ws.Cells.Item(48, 1).Value = 45889.48958333334
$ws.Cells.Item(48, 2).Value = 2301
$ws.Cells.Item(49, 1).Value = 45889.5
$ws.Cells.Item(49, 2).Value = 2300
$ws.Cells.Item(50, 1).Value = 45889.51041666666
$ws.Cells.Item(50, 2).Value = 2275
$ws.Cells.Item(51, 1).Value = 45889.52083333334
$ws.Cells.Item(51, 2).Value = 2263
$ws.Cells.Item(52, 1).Value = 45889.53125
$ws.Cells.Item(52, 2).Value = 2261
$ws.Cells.Item(53, 1).Value = 45889.54166666666
$ws.Cells.Item(53, 2).Value = 2252
$ws.Cells.Item(54, 1).Value = 45889.55208333334
$ws.Cells.Item(54, 2).Value = 2189
$ws.Cells.Item(55, 1).Value = 45889.5625
$ws.Cells.Item(55, 2).Value = 2181
$ws.Cells.Item(56, 1).Value = 45889.57291666666
$ws.Cells.Item(56, 2).Value = 2171
$ws.Cells.Item(57, 1).Value = 45889.58333333334
$ws.Cells.Item(57, 2).Value = 2157
$ws.Cells.Item(58, 1).Value = 45889.59375
$ws.Cells.Item(58, 2).Value = 2010
$ws.Cells.Item(59, 1).Value = 45889.60416666666
$ws.Cells.Item(59, 2).Value = 1990
$ws.Cells.Item(60, 1).Value = 45889.61458333334
$ws.Cells.Item(60, 2).Value = 1971
$ws.Cells.Item(61, 1).Value = 45889.625
$ws.Cells.Item(61, 2).Value = 1949
$ws.Cells.Item(62, 1).Value = 45889.63541666666
$ws.Cells.Item(62, 2).Value = 1716
$ws.Cells.Item(63, 1).Value = 45889.64583333334
$ws.Cells.Item(63, 2).Value = 1679
$ws.Cells.Item(64, 1).Value = 45889.65625
$ws.Cells.Item(64, 2).Value = 1636
$ws.Cells.Item(65, 1).Value = 45889.66666666666
$ws.Cells.Item(65, 2).Value = 1601
$ws.Cells.Item(66, 1).Value = 45889.67708333334
$ws.Cells.Item(66, 2).Value = 1246
$ws.Cells.Item(67, 1).Value = 45889.6875
$ws.Cells.Item(67, 2).Value = 1200
$ws.Cells.Item(68, 1).Value = 45889.69791666666
$ws.Cells.Item(68, 2).Value = 1134
$ws.Cells.Item(69, 1).Value = 45889.70833333334
$ws.Cells.Item(69, 2).Value = 1091
$ws.Cells.Item(70, 1).Value = 45889.71875
$ws.Cells.Item(70, 2).Value = 614
$ws.Cells.Item(71, 1).Value = 45889.72916666666
$ws.Cells.Item(71, 2).Value = 575
$ws.Cells.Item(72, 1).Value = 45889.73958333334
$ws.Cells.Item(72, 2).Value = 512
$ws.Cells.Item(73, 1).Value = 45889.75
$ws.Cells.Item(73, 2).Value = 480
$ws.Cells.Item(74, 1).Value = 45889.76041666666
$ws.Cells.Item(74, 2).Value = 183
$ws.Cells.Item(75, 1).Value = 45889.77083333334
$ws.Cells.Item(75, 2).Value = 152
$ws.Cells.Item(76, 1).Value = 45889.78125
$ws.Cells.Item(76, 2).Value = 122
$ws.Cells.Item(77, 1).Value = 45889.79166666666
$ws.Cells.Item(77, 2).Value = 109
$ws.Cells.Item(78, 1).Value = 45889.80208333334
$ws.Cells.Item(78, 2).Value = 21
$ws.Cells.Item(79, 1).Value = 45889.8125
$ws.Cells.Item(79, 2).Value = 18
$ws.Cells.Item(80, 1).Value = 45889.82291666666
$ws.Cells.Item(80, 2).Value = 18
$ws.Cells.Item(81, 1).Value = 45889.83333333334
$ws.Cells.Item(81, 2).Value = 18
$ws.Cells.Item(82, 1).Value = 45889.84375
$ws.Cells.Item(82, 2).Value = 7
$ws.Cells.Item(83, 1).Value = 45889.85416666666
$ws.Cells.Item(83, 2).Value = 7
$ws.Cells.Item(84, 1).Value = 45889.86458333334
$ws.Cells.Item(84, 2).Value = 7
$ws.Cells.Item(85, 1).Value = 45889.875
$ws.Cells.Item(85, 2).Value = 7
$ws.Cells.Item(86, 1).Value = 45889.88541666666
$ws.Cells.Item(86, 2).Value = 1
$ws.Cells.Item(87, 1).Value = 45889.89583333334
$ws.Cells.Item(87, 2).Value = 1
$ws.Cells.Item(88, 1).Value = 45889.90625
$ws.Cells.Item(88, 2).Value = 1
$ws.Cells.Item(89, 1).Value = 45889.91666666666
$ws.Cells.Item(89, 2).Value = 1
$ws.Cells.Item(90, 1).Value = 45889.92708333334
$ws.Cells.Item(91, 1).Value = 45889.9375
$ws.Cells.Item(92, 1).Value = 45889.94791666666
$ws.Cells.Item(93, 1).Value = 45889.95833333334
$ws.Cells.Item(94, 1).Value = 45889.96875
$ws.Cells.Item(95, 1).Value = 45889.97916666666
$ws.Cells.Item(96, 1).Value = 45889.98958333334
$ws.Cells.Item(97, 1).Value = 45890
